$d = $word.ActiveDocument

# The title paragraph's first run contains the text "Oct" (no trailing
# space - the following run already starts with a leading space before
# "Docker part"). Remove that run entirely so the paragraph begins with
# " Docker part6 HW ".
$para = $d.Paragraphs(1)
$run = $para.Range.Words(1)

$d.Content.Find.Execute("Oct", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng = $d.Content.Duplicate
$rng.Start = $d.Content.Start
$rng.End = $d.Content.Start + 3
$rng.Delete()
